$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.74649233167628
$ws.Range("C2").Value = 5.280149767803496
$ws.Range("D2").Value = 6.404481587758474
$ws.Range("F2").Value = 32.22342665823589
$ws.Range("G2").Value = 3.678521223139385
$ws.Range("K2").Value = 12.14800977874676
$ws.Range("N2").Value = 20.27594957129025

$ws.Range("B3").Value = 12.50336651230864
$ws.Range("C3").Value = 5.080188451177959
$ws.Range("D3").Value = 6.402417816100481
$ws.Range("F3").Value = 32.19741226299392
$ws.Range("G3").Value = 3.681109996171884
$ws.Range("K3").Value = 11.98651058806879
$ws.Range("N3").Value = 20.34288248931044

$ws.Range("B4").Value = 12.35525243798511
$ws.Range("C4").Value = 4.95477656679498
$ws.Range("D4").Value = 6.401172251513683
$ws.Range("F4").Value = 32.19004152305256
$ws.Range("G4").Value = 3.682782756929701
$ws.Range("K4").Value = 11.88942801688785
$ws.Range("N4").Value = 20.38596737323946

$ws.Range("B5").Value = 12.29527813414207
$ws.Range("C5").Value = 4.903096272069977
$ws.Range("D5").Value = 6.400670362410817
$ws.Range("F5").Value = 32.18920083301195
$ws.Range("G5").Value = 3.683485423341819
$ws.Range("K5").Value = 11.85043798169526
$ws.Range("N5").Value = 20.40402575848976

$ws.Range("B6").Value = 12.28534524851131
$ws.Range("C6").Value = 4.894482892634818
$ws.Range("D6").Value = 6.400587375285747
$ws.Range("F6").Value = 32.18919183533412
$ws.Range("G6").Value = 3.683603371156198
$ws.Range("K6").Value = 11.84399967465767
$ws.Range("N6").Value = 20.40705462784966

$ws.Range("B7").Value = 12.35444193168133
$ws.Range("C7").Value = 4.954081788676127
$ws.Range("D7").Value = 6.40116545950067
$ws.Range("F7").Value = 32.19002142908585
$ws.Range("G7").Value = 3.682792148201606
$ws.Range("K7").Value = 11.88889980381079
$ws.Range("N7").Value = 20.38620888519992

$ws.Range("B8").Value = 12.66247433428772
$ws.Range("C8").Value = 5.211800470866508
$ws.Range("D8").Value = 6.403765668723506
$ws.Range("F8").Value = 32.21267192782567
$ws.Range("G8").Value = 3.679396597395532
$ws.Range("K8").Value = 12.09192328107408
$ws.Range("N8").Value = 20.29861613568332

$ws.Range("B9").Value = 13.27187467136383
$ws.Range("C9").Value = 5.692821007805753
$ws.Range("D9").Value = 6.409025171864348
$ws.Range("F9").Value = 32.32527725540884
$ws.Range("G9").Value = 3.673395213336245
$ws.Range("K9").Value = 12.50434025881685
$ws.Range("N9").Value = 20.14256821035945

$ws.Range("B10").Value = 13.7177758810058
$ws.Range("C10").Value = 6.027328777404513
$ws.Range("D10").Value = 6.41297315659182
$ws.Range("F10").Value = 32.44938002435494
$ws.Range("G10").Value = 3.669382181780044
$ws.Range("K10").Value = 12.81307486890999
$ws.Range("N10").Value = 20.03743077773014

$ws.Range("B11").Value = 13.91921171054963
$ws.Range("C11").Value = 6.174704356893296
$ws.Range("D11").Value = 6.414784519428197
$ws.Range("F11").Value = 32.51474137336991
$ws.Range("G11").Value = 3.667641616755994
$ws.Range("K11").Value = 12.95414167053173
$ws.Range("N11").Value = 19.99165082407817

$ws.Range("B12").Value = 13.99520792756699
$ws.Range("C12").Value = 6.229773830997799
$ws.Range("D12").Value = 6.415472408404242
$ws.Range("F12").Value = 32.54076294096668
$ws.Range("G12").Value = 3.666994657452095
$ws.Range("K12").Value = 13.00759783976138
$ws.Range("N12").Value = 19.97460849664718

$ws.Range("B13").Value = 13.97885461815522
$ws.Range("C13").Value = 6.21794719179912
$ws.Range("D13").Value = 6.415324175737188
$ws.Range("F13").Value = 32.53510239838619
$ws.Range("G13").Value = 3.667133452161694
$ws.Range("K13").Value = 12.99608425934105
$ws.Range("N13").Value = 19.97826582519238

$ws.Range("B14").Value = 13.92547008800935
$ws.Range("C14").Value = 6.179250014073987
$ws.Range("D14").Value = 6.414841074501785
$ws.Range("F14").Value = 32.5168567820085
$ws.Range("G14").Value = 3.667588147796207
$ws.Range("K14").Value = 12.95853909018256
$ws.Range("N14").Value = 19.99024286346455

$ws.Range("B15").Value = 13.89273131883484
$ws.Range("C15").Value = 6.155449358076529
$ws.Range("D15").Value = 6.41454540961732
$ws.Range("F15").Value = 32.50584596020614
$ws.Range("G15").Value = 3.667868242989506
$ws.Range("K15").Value = 12.93554491030173
$ws.Range("N15").Value = 19.99761734759533

$ws.Range("B16").Value = 13.70457625850656
$ws.Range("C16").Value = 6.017596626442176
$ws.Range("D16").Value = 6.412855068289232
$ws.Range("F16").Value = 32.44528690516127
$ws.Range("G16").Value = 3.669497636260655
$ws.Range("K16").Value = 12.80386384745011
$ws.Range("N16").Value = 20.04046373128578

$ws.Range("B17").Value = 13.58873004119262
$ws.Range("C17").Value = 5.931763612969654
$ws.Range("D17").Value = 6.411821865024062
$ws.Range("F17").Value = 32.41041042795882
$ws.Range("G17").Value = 3.670518935681863
$ws.Range("K17").Value = 12.72320356711413
$ws.Range("N17").Value = 20.06727239502406

$ws.Range("B18").Value = 13.52197052072903
$ws.Range("C18").Value = 5.881946030736785
$ws.Range("D18").Value = 6.411229058443679
$ws.Range("F18").Value = 32.39118943090179
$ws.Range("G18").Value = 3.671114363004392
$ws.Range("K18").Value = 12.67687150117901
$ws.Range("N18").Value = 20.08288481042872

$ws.Range("B19").Value = 13.49934746966032
$ws.Range("C19").Value = 5.865003192730474
$ws.Range("D19").Value = 6.411028603568289
$ws.Range("F19").Value = 32.38482591327666
$ws.Range("G19").Value = 3.671317340978116
$ws.Range("K19").Value = 12.66119644755888
$ws.Range("N19").Value = 20.08820404396398

$ws.Range("B20").Value = 13.60107589883697
$ws.Range("C20").Value = 5.940947499400581
$ws.Range("D20").Value = 6.411931701503584
$ws.Range("F20").Value = 32.41403631922446
$ws.Range("G20").Value = 3.670409388820332
$ws.Range("K20").Value = 12.73178399869022
$ws.Range("N20").Value = 20.06439862138877

$ws.Range("B21").Value = 13.9411587267265
$ws.Range("C21").Value = 6.190636714474242
$ws.Range("D21").Value = 6.414982921682173
$ws.Range("F21").Value = 32.52218156625808
$ws.Range("G21").Value = 3.667454263306846
$ws.Range("K21").Value = 12.96956641591261
$ws.Range("N21").Value = 19.98671695722355

$ws.Range("B22").Value = 14.1617375932124
$ws.Range("C22").Value = 6.349500368377536
$ws.Range("D22").Value = 6.41698841188138
$ws.Range("F22").Value = 32.60026107084333
$ws.Range("G22").Value = 3.665593735869086
$ws.Range("K22").Value = 13.1251646994357
$ws.Range("N22").Value = 19.93765811799199

$ws.Range("B23").Value = 14.04419027185465
$ws.Range("C23").Value = 6.265121969469338
$ws.Range("D23").Value = 6.415917087859077
$ws.Range("F23").Value = 32.55791523598428
$ws.Range("G23").Value = 3.666580276087686
$ws.Range("K23").Value = 13.04211815714677
$ws.Range("N23").Value = 19.963685512129

$ws.Range("B24").Value = 13.59549482246432
$ws.Range("C24").Value = 5.936796927746651
$ws.Range("D24").Value = 6.411882040660441
$ws.Range("F24").Value = 32.41239446769249
$ws.Range("G24").Value = 3.670458889194143
$ws.Range("K24").Value = 12.72790465508316
$ws.Range("N24").Value = 20.06569723208911

$ws.Range("B25").Value = 13.10697851477175
$ws.Range("C25").Value = 5.565749544275298
$ws.Range("D25").Value = 6.407586618059033
$ws.Range("F25").Value = 32.2875277179829
$ws.Range("G25").Value = 3.674948848218013
$ws.Range("K25").Value = 12.39154843797377
$ws.Range("N25").Value = 20.18310787214158
